$d = $word.ActiveDocument
$sec = $d.Sections(1)
$hdr = $sec.Headers(2)
for ($k = 1; $k -le $hdr.Shapes.Count; $k++) {
    $shp = $hdr.Shapes.Item($k)
    Write-Host "Shape $k :" $shp.Name
    try {
        $t = $shp.TextFrame.HasText
        Write-Host "  HasText:" $t
    } catch {
        Write-Host "  ERROR HasText:" $_
    }
}

Write-Host "--- header1 ---"
$hdr1 = $sec.Headers(1)
for ($k = 1; $k -le $hdr1.Shapes.Count; $k++) {
    $shp = $hdr1.Shapes.Item($k)
    Write-Host "Shape $k :" $shp.Name
}

Write-Host "--- footer1 (sec.Footers(1)) ---"
$ftr1 = $sec.Footers(1)
for ($k = 1; $k -le $ftr1.Shapes.Count; $k++) {
    $shp = $ftr1.Shapes.Item($k)
    Write-Host "Shape $k :" $shp.Name
}
